$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2441.111
$ws.Range("I15").Value = 2441.111
$ws.Range("K15").Value = 7323.333
$ws.Range("M15").Value = -7154.333
$ws.Range("H17").Value = 1668960.9
$ws.Range("J17").Value = 1697237.1
$ws.Range("L17").Value = 5091711.300000001
$ws.Range("N17").Value = -5092047.300000001
$ws.Range("H110").Value = 45700
$ws.Range("J110").Value = 45700
$ws.Range("L110").Value = 45700
$ws.Range("N110").Value = -53880
$ws.Range("H116").Value = 5713.25
$ws.Range("I116").Value = 2795
$ws.Range("J116").Value = 6130.143
$ws.Range("K116").Value = 2795
$ws.Range("L116").Value = 6130.143
$ws.Range("M116").Value = 647
$ws.Range("N116").Value = -13014.143
$ws.Range("H129").Value = 189486.12
$ws.Range("J129").Value = 200833.3
$ws.Range("L129").Value = 602499.8999999999
$ws.Range("N129").Value = -612499.8999999999
$ws.Range("H137").Value = 65475.125
$ws.Range("I137").Value = 3570.1
$ws.Range("K137").Value = 10710.3
$ws.Range("M137").Value = -8160.299999999999
$ws.Range("H138").Value = 1913.475
$ws.Range("I138").Value = 1275.5161
$ws.Range("J138").Value = 2317.0815
$ws.Range("K138").Value = 3826.5483
$ws.Range("L138").Value = 6951.244499999999
$ws.Range("M138").Value = 1313.4517
$ws.Range("N138").Value = -17231.2445

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1218.8438
$ws.Range("I2").Value = 1175.762
$ws.Range("J2").Value = 1301.091
$ws.Range("K2").Value = 1175.762
$ws.Range("L2").Value = 1301.091
$ws.Range("M2").Value = -1062.762
$ws.Range("N2").Value = -1527.091
$ws.Range("H32").Value = 29808.605
$ws.Range("I32").Value = 30436.473
$ws.Range("K32").Value = 30436.473
$ws.Range("M32").Value = -30149.473
$ws.Range("H34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("H45").Value = 2086.04
$ws.Range("I45").Value = 2202.6428
$ws.Range("K45").Value = 2202.6428
$ws.Range("M45").Value = -1825.6428
$ws.Range("H61").Value = 2838.318
$ws.Range("I61").Value = 2273.1177
$ws.Range("K61").Value = 2273.1177
$ws.Range("M61").Value = -2061.1177
$ws.Range("H116").Value = 1218.8438
$ws.Range("I116").Value = 1175.762
$ws.Range("J116").Value = 1301.091
$ws.Range("K116").Value = 1175.762
$ws.Range("L116").Value = 1301.091
$ws.Range("M116").Value = 1118.238
$ws.Range("N116").Value = -5889.091
$ws.Range("H136").Value = 2838.318
$ws.Range("I136").Value = 2273.1177
$ws.Range("K136").Value = 6819.353099999999
$ws.Range("M136").Value = -4269.353099999999
$ws.Range("N34").ClearContents()

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1218.8438
$ws.Range("I3").Value = 1175.762
$ws.Range("J3").Value = 1301.091
$ws.Range("K3").Value = 1175.762
$ws.Range("L3").Value = 1301.091
$ws.Range("M3").Value = -1061.762
$ws.Range("N3").Value = -1529.091
$ws.Range("H94").Value = 1341.5116
$ws.Range("I94").Value = 1099.359
$ws.Range("J94").Value = 3702.5
$ws.Range("K94").Value = 1099.359
$ws.Range("L94").Value = 3702.5
$ws.Range("M94").Value = -648.3589999999999
$ws.Range("N94").Value = -4604.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 216.66667
$ws.Range("J22").Value = 300
$ws.Range("L22").Value = 300
$ws.Range("N22").Value = -1000
$ws.Range("H31").Value = 13055.929
$ws.Range("I31").Value = 16892.947
$ws.Range("K31").Value = 16892.947
$ws.Range("M31").Value = -16597.947
$ws.Range("H34").Value = 13055.929
$ws.Range("I34").Value = 16892.947
$ws.Range("K34").Value = 16892.947
$ws.Range("M34").Value = -16690.947
$ws.Range("H35").Value = 3500
$ws.Range("I35").Value = 3500
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 3500
$ws.Range("L35").Value = 0
$ws.Range("H58").Value = 14532
$ws.Range("I58").Value = 1040.4814
$ws.Range("K58").Value = 1040.4814
$ws.Range("M58").Value = -837.4813999999999
$ws.Range("H94").Value = 2338.5908
$ws.Range("J94").Value = 3381.9
$ws.Range("L94").Value = 3381.9
$ws.Range("N94").Value = -4283.9
$ws.Range("H122").Value = 1840.4615
$ws.Range("I122").Value = 2272
$ws.Range("J122").Value = 1150
$ws.Range("K122").Value = 6816
$ws.Range("L122").Value = 3450
$ws.Range("M122").Value = -4366
$ws.Range("N122").Value = -8350
$ws.Range("H136").Value = 14532
$ws.Range("I136").Value = 1040.4814
$ws.Range("K136").Value = 3121.4442
$ws.Range("M136").Value = -571.4441999999999
$ws.Range("N35").ClearContents()

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 715.26
$ws.Range("J131").Value = 744.25275
$ws.Range("L131").Value = 2232.75825
$ws.Range("N131").Value = -12312.75825
$ws.Range("H138").Value = 112488.7
$ws.Range("J138").Value = 177713.83
$ws.Range("L138").Value = 533141.49
$ws.Range("N138").Value = -543421.49

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 298.5
$ws.Range("I97").Value = 307.7647
$ws.Range("J97").Value = 267
$ws.Range("K97").Value = 307.7647
$ws.Range("L97").Value = 267
$ws.Range("M97").Value = 188.2353
$ws.Range("N97").Value = -1259
$ws.Range("H132").Value = 61971.81
$ws.Range("I132").Value = 54495.45
$ws.Range("J132").Value = 86893
$ws.Range("K132").Value = 163486.35
$ws.Range("L132").Value = 260679
$ws.Range("M132").Value = -160956.35
$ws.Range("N132").Value = -265739

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7340.2666
$ws.Range("I7").Value = 7945
$ws.Range("J7").Value = 5677.25
$ws.Range("K7").Value = 7945
$ws.Range("L7").Value = 5677.25
$ws.Range("M7").Value = -7833
$ws.Range("N7").Value = -5901.25
$ws.Range("H126").Value = 7340.2666
$ws.Range("I126").Value = 7945
$ws.Range("J126").Value = 5677.25
$ws.Range("K126").Value = 23835
$ws.Range("L126").Value = 17031.75
$ws.Range("M126").Value = -21365
$ws.Range("N126").Value = -21971.75
$ws.Range("H132").Value = 2742.1052
$ws.Range("I132").Value = 1864.4286
$ws.Range("K132").Value = 5593.2858
$ws.Range("M132").Value = -3063.2858

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1798.9333
$ws.Range("I132").Value = 978.875
$ws.Range("J132").Value = 2736.1428
$ws.Range("K132").Value = 2936.625
$ws.Range("L132").Value = 8208.428400000001
$ws.Range("M132").Value = -406.625
$ws.Range("N132").Value = -13268.4284
$ws.Range("H136").Value = 28573400
$ws.Range("I136").Value = 45456370
$ws.Range("K136").Value = 136369110
$ws.Range("M136").Value = -136366560
